$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.428.07"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "1.796.04"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  +0.15%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.85"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("E7").Value = "  +0.10%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.45"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  +1.60%  "
$ws.Range("E9").Value = "  +1.21%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0694"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +0.44%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0949"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "2.056.03"
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.08"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").Value = "1.793.03"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").Value = "34.388.37"
$ws.Range("E16").Value = "  +1.05%  "
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.23"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  +1.15%  "
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.37"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").Value = "0.0₃0803"
$ws.Range("E19").Value = "  +3.12%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "246.92"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +0.68%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.04"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +1.71%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.15"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +1.17%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  +1.37%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.68"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +0.81%  "
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.20"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +0.81%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.40"
$ws.Range("D27").Style = $origStyle
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  +0.39%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0522"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +1.07%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.77"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +3.16%  "
$ws.Range("E33").Value = "  +7.73%  "
$ws.Range("E34").Value = "  +1.48%  "
$ws.Range("D35").Value = "1.445.42"
$ws.Range("E35").Value = "  -0.61%  "
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.62"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  +8.20%  "
$ws.Range("E37").Value = "  +3.31%  "
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("E39").Value = "  -0.58%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "84.09"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +4.86%  "
$ws.Range("E41").Value = "  +1.41%  "
$ws.Range("E42").Value = "  +3.12%  "
$ws.Range("E43").Value = "  +1.85%  "
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.84"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("E45").Value = "  +3.68%  "
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "1.951.48"
$ws.Range("E48").Value = "  +0.30%  "
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.50"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("D51").Value = "0.0₆0129"
$ws.Range("E51").Value = "  -4.64%  "
